$p = $ppt.ActivePresentation
$layout = $p.SlideMaster.CustomLayouts.Item(2)

# --- Slide 2 : Fink-dataset-monitor (enjeux / outils) ---
$s2 = $p.Slides.AddSlide(2, $layout)
$s2Title = $s2.Shapes.Item(1)
$s2Body = $s2.Shapes.Item(2)
$s2Title.TextFrame.TextRange.Text = "Fink-"
$s2Title.TextFrame.TextRange.LanguageID = "fr-FR"
$__tmp = $s2Title.TextFrame.TextRange.InsertAfter("dataset")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s2Title.TextFrame.TextRange.InsertAfter("-monitor")
$__tmp.LanguageID = "fr-FR"
$s2Body.TextFrame.TextRange.Text = "Plusieurs enjeux"
$s2Body.TextFrame.TextRange.LanguageID = "fr-FR"
$__tmp = $s2Body.TextFrame.TextRange.InsertAfter([char]13 + "Définir une architecture client-serveur pour donner accès aux données du broker pour les utilisateurs non-locaux")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s2Body.TextFrame.TextRange.InsertAfter([char]13 + "Ne pas obliger les utilisateurs à posséder un compte sur le (futur) cluster qui accueil le broker")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s2Body.TextFrame.TextRange.InsertAfter([char]13 + "Structure ouverte (web ")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s2Body.TextFrame.TextRange.InsertAfter("based")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s2Body.TextFrame.TextRange.InsertAfter(") pour augmenter les outils de présentation des données")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s2Body.TextFrame.TextRange.InsertAfter([char]13 + "Outils:")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s2Body.TextFrame.TextRange.InsertAfter([char]13 + "Livy")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s2Body.TextFrame.TextRange.InsertAfter([char]13 + "Hbase")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s2Body.TextFrame.TextRange.InsertAfter([char]13 + "Authentification")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s2Body.TextFrame.TextRange.InsertAfter([char]13 + "Séparer le monde du cluster du monde des utilisateurs")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s2Body.TextFrame.TextRange.InsertAfter([char]13 + "Accès au cluster à travers une Gateway")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s2Body.TextFrame.TextRange.InsertAfter([char]13 + "Construire une organisation autonome (gestion via ")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s2Body.TextFrame.TextRange.InsertAfter("Hbase")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s2Body.TextFrame.TextRange.InsertAfter(")")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s2Body.TextFrame.TextRange.InsertAfter([char]13 + "Transmission des identifiants des utilisateurs via un protocole ad-hoc")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s2Body.TextFrame.TextRange.InsertAfter([char]13 + "Traçabilité")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s2Body.TextFrame.TextRange.InsertAfter([char]13 + "Accounting")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s2Body.TextFrame.TextRange.InsertAfter(" (ressources)")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s2Body.TextFrame.TextRange.InsertAfter([char]13 + "")
$__tmp.LanguageID = "fr-FR"
$s2Body.TextFrame.TextRange.Paragraphs(2, 1).IndentLevel = 2
$s2Body.TextFrame.TextRange.Paragraphs(3, 1).IndentLevel = 2
$s2Body.TextFrame.TextRange.Paragraphs(4, 1).IndentLevel = 2
$s2Body.TextFrame.TextRange.Paragraphs(6, 1).IndentLevel = 2
$s2Body.TextFrame.TextRange.Paragraphs(7, 1).IndentLevel = 2
$s2Body.TextFrame.TextRange.Paragraphs(9, 1).IndentLevel = 2
$s2Body.TextFrame.TextRange.Paragraphs(10, 1).IndentLevel = 3
$s2Body.TextFrame.TextRange.Paragraphs(11, 1).IndentLevel = 3
$s2Body.TextFrame.TextRange.Paragraphs(12, 1).IndentLevel = 3
$s2Body.TextFrame.TextRange.Paragraphs(13, 1).IndentLevel = 4
$s2Body.TextFrame.TextRange.Paragraphs(14, 1).IndentLevel = 4
$s2Body.TextFrame.TextRange.Paragraphs(15, 1).IndentLevel = 3
$s2Body.TextFrame.AutoSize = 2

# --- Slide 3 : Fink-dataset-monitor (principe / suite) ---
$s3 = $p.Slides.AddSlide(3, $layout)
$s3Title = $s3.Shapes.Item(1)
$s3Body = $s3.Shapes.Item(2)
$s3Title.TextFrame.TextRange.Text = "Fink-"
$s3Title.TextFrame.TextRange.LanguageID = "fr-FR"
$__tmp = $s3Title.TextFrame.TextRange.InsertAfter("dataset")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Title.TextFrame.TextRange.InsertAfter("-monitor")
$__tmp.LanguageID = "fr-FR"
$s3Body.TextFrame.TextRange.Text = "Principe"
$s3Body.TextFrame.TextRange.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter([char]13 + "Gestion d’une organisation d’autorisations d’identification (")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter("Hbase")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter(")")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter([char]13 + "Accès aux données via ")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter("Livy")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter("/")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter("Hbase")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter("(")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter("dataframes")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter(")")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter([char]13 + "Gestion des sessions ")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter("Livy")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter("/")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter("Yarn")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter(" maintenues ")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter([char]13 + "Un utilisateur peut relancer son client et se reconnecter à sa session ouverte")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter([char]13 + "Les administrateurs peuvent associer les utilisateurs aux sessions ")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter("Livy")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter(" pour la gestion des ressources allouées")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter([char]13 + "Un prototype opérationnel existe qui implémente le modèle ")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter([char]13 + "Suite:")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter([char]13 + "Besoin d’expertise pour définir les outils d’authentification (LDAP? Xxx?)")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter([char]13 + "Développements sur la structure «")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter("pluggin")
$__tmp.LanguageID = "fr-FR"
$__tmp = $s3Body.TextFrame.TextRange.InsertAfter("» du client web")
$__tmp.LanguageID = "fr-FR"
$s3Body.TextFrame.TextRange.Paragraphs(2, 1).IndentLevel = 2
$s3Body.TextFrame.TextRange.Paragraphs(3, 1).IndentLevel = 2
$s3Body.TextFrame.TextRange.Paragraphs(4, 1).IndentLevel = 2
$s3Body.TextFrame.TextRange.Paragraphs(5, 1).IndentLevel = 3
$s3Body.TextFrame.TextRange.Paragraphs(6, 1).IndentLevel = 3
$s3Body.TextFrame.TextRange.Paragraphs(9, 1).IndentLevel = 2
$s3Body.TextFrame.TextRange.Paragraphs(10, 1).IndentLevel = 2
$s3Body.TextFrame.AutoSize = 2

